$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.133.07"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "1.637.18"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "216.75"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  +2.31%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +0.28%  "

$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Value = "1.865.47"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "1.632.96"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "'0.540"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").Value = "66.58"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").Value = "27.124.38"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").Value = "216.59"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  +3.49%  "

$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").Value = "146.89"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  +1.84%  "

$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.37%  "

$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "1.298.13"
$ws.Range("E34").Value = "  +2.65%  "

$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("E37").Value = "  -0.55%  "

$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +1.54%  "

$ws.Range("E39").Value = "  +1.89%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +6.01%  "

$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "1.775.62"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").Value = "61.57"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").Value = "91.38"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("E51").Value = "  -0.41%  "
